# Update: refresh the "Förändrad" (changed) date stamp for every existing
# data row (2-407) from 2023-09-23 (45192) to 2023-10-03 (45202), and
# append three new felling-notification rows (408-410) for VINGÅKER.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump column C ("Förändrad") for all existing data rows (2..407) ---
$ws.Range("C2:C407").Value2 = 45202

# Row 407 previously had the implicit default row height; the source file
# now pins it explicitly (ht="15" customHeight="1"), matching rows 408/409
# added below.
$ws.Rows.Item(407).RowHeight = 15

# --- 2. Append the three new rows ---

# Row 408
$ws.Cells.Item(408, 1).Value2 = "A 45881-2023"
$ws.Cells.Item(408, 2).Value2 = 45195
$ws.Cells.Item(408, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(408, 3).Value2 = 45202
$ws.Cells.Item(408, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(408, 4).Value2 = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(408, 5).Value2 = "VINGÅKER"
$ws.Cells.Item(408, 7).Value2 = 3.8
$ws.Cells.Item(408, 8).Value2 = 0
$ws.Cells.Item(408, 9).Value2 = 0
$ws.Cells.Item(408, 10).Value2 = 0
$ws.Cells.Item(408, 11).Value2 = 0
$ws.Cells.Item(408, 12).Value2 = 0
$ws.Cells.Item(408, 13).Value2 = 0
$ws.Cells.Item(408, 14).Value2 = 0
$ws.Cells.Item(408, 15).Value2 = 0
$ws.Cells.Item(408, 16).Value2 = 0
$ws.Cells.Item(408, 17).Value2 = 0
$ws.Cells.Item(408, 18).WrapText = $true
$ws.Rows.Item(408).RowHeight = 15

# Row 409
$ws.Cells.Item(409, 1).Value2 = "A 46159-2023"
$ws.Cells.Item(409, 2).Value2 = 45196
$ws.Cells.Item(409, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(409, 3).Value2 = 45202
$ws.Cells.Item(409, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(409, 4).Value2 = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(409, 5).Value2 = "VINGÅKER"
$ws.Cells.Item(409, 7).Value2 = 9.3
$ws.Cells.Item(409, 8).Value2 = 0
$ws.Cells.Item(409, 9).Value2 = 0
$ws.Cells.Item(409, 10).Value2 = 0
$ws.Cells.Item(409, 11).Value2 = 0
$ws.Cells.Item(409, 12).Value2 = 0
$ws.Cells.Item(409, 13).Value2 = 0
$ws.Cells.Item(409, 14).Value2 = 0
$ws.Cells.Item(409, 15).Value2 = 0
$ws.Cells.Item(409, 16).Value2 = 0
$ws.Cells.Item(409, 17).Value2 = 0
$ws.Cells.Item(409, 18).WrapText = $true
$ws.Rows.Item(409).RowHeight = 15

# Row 410
$ws.Cells.Item(410, 1).Value2 = "A 46052-2023"
$ws.Cells.Item(410, 2).Value2 = 45196
$ws.Cells.Item(410, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(410, 3).Value2 = 45202
$ws.Cells.Item(410, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(410, 4).Value2 = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(410, 5).Value2 = "VINGÅKER"
$ws.Cells.Item(410, 7).Value2 = 3.3
$ws.Cells.Item(410, 8).Value2 = 0
$ws.Cells.Item(410, 9).Value2 = 0
$ws.Cells.Item(410, 10).Value2 = 0
$ws.Cells.Item(410, 11).Value2 = 0
$ws.Cells.Item(410, 12).Value2 = 0
$ws.Cells.Item(410, 13).Value2 = 0
$ws.Cells.Item(410, 14).Value2 = 0
$ws.Cells.Item(410, 15).Value2 = 0
$ws.Cells.Item(410, 16).Value2 = 0
$ws.Cells.Item(410, 17).Value2 = 0
$ws.Cells.Item(410, 18).WrapText = $true

Write-Host "Updated C2:C407 and appended rows 408-410"
